# "worked out the multiple member print, need to fix the dates of the week to align"
#
# Duplicate the "Sheet4" worksheet (Excel's default "Move or Copy -> Create a
# copy" placed immediately before the source sheet, which is how Excel names
# it "Sheet4 (2)" and inserts it ahead of "Sheet4"). Then trim the new copy
# down to just the header + first member's rows (rows 1-5 of 8), since this
# is meant to be the per-member print-out sheet and the extra member rows
# from the template don't belong there yet.

$wb = $excel.ActiveWorkbook

# Locate the source sheet by name (index-based handles shift once sheets are
# inserted, so re-resolve everything we need up front).
$source = $wb.Worksheets.Item("Sheet4")

# Copy "Sheet4" to just before itself -> Excel places the new sheet first
# and names it "Sheet4 (2)".
$source.Copy($source)

# Re-fetch handles by name: after the copy, any previously held reference to
# "Sheet4" now actually addresses whichever sheet sits at that old index
# (i.e. the freshly-inserted copy), so grab both sheets fresh by name.
$copy = $wb.Worksheets.Item("Sheet4 (2)")
$original = $wb.Worksheets.Item("Sheet4")

# The template copy only needs the header row plus the first member's row
# (rows 1-5); drop the remaining member rows (6-8) that were duplicated in.
$copy.Rows("6:8").Delete()

# Restore sensible per-sheet selections: the original "Sheet4" keeps working
# at J16, while the new copy is left active with D16 selected.
$original.Activate()
$original.Range("J16").Select()

$copy.Activate()
$copy.Range("D16").Select()
